# "ajuste na url de acesso a aplicacao"
#
# The README shows the app's access URL as:
#   http://localhost:8080/enderecamentoapp/
# This changes the last path segment from "enderecamentoapp" to "webapp",
# rendered as two underlined runs ("web" + "app") so the change stands out
# visually in the document, matching the rest of that sentence's styling.

$d = $word.ActiveDocument

# Locate the paragraph that contains the access URL (there are other
# paragraphs mentioning "enderecamentoapp" - e.g. the deploy/undeploy/test
# instructions - that must stay untouched).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*http://localhost:8080/enderecamentoapp/*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Narrow a Range down to just the "enderecamentoapp" word inside that
    # paragraph.
    $r = $target.Range.Duplicate
    $found = $r.Find.Execute("enderecamentoapp", $false, $false, $false, `
                              $false, $false, $true, 1, $false, "", 0)

    if ($found) {
        $start = $r.Start

        # Swap the text for "webapp" ...
        $r.Text = "webapp"

        # ... then underline it as two separate runs, "web" and "app".
        $rWeb = $d.Range($start, $start + 3)
        $rWeb.Font.Underline = 1

        $rApp = $d.Range($start + 3, $start + 6)
        $rApp.Font.Underline = 1
    }
}

Write-Output "Done: $($target.Range.Text)"
